$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Range("H1").Value = "NRMSE"

# Match column width of other data columns (stored width 12.7109375 in the
# original workbook XML; closest achievable value through the COM
# ColumnWidth property, which snaps to whole-pixel increments)
$ws.Columns.Item(8).ColumnWidth = 11.83

# Updated data values for columns A-G, rows 2-6
$data = @(
    @(3.2895343668952766, 2.3127057387057386, 0.4520745120813951, 0.67236486529368489, 0.64564649321768253, 0.54792548791860485, 0.745739696638961, 0.23941298157898669),
    @(2.9204697635116927, 2.1854963624338621, 0.27616617092892587, 0.52551514814411004, 0.51678797882096528, 0.72383382907107419, 0.85599777326448456, 0.20380110003570778),
    @(3.70846076494434, 2.76842871189774, 0.44457278403895834, 0.66676291441483027, 0.62156010594920064, 0.55542721596104161, 0.75000970226409058, 0.25261994311609942),
    @(3.9495173561721337, 3.1118916790674596, 0.47903875079228225, 0.69212625350602197, 0.72260343180482045, 0.5209612492077178, 0.72319320180540436, 0.28200766556030948),
    @(3.4790983134161579, 2.6591716865930883, 0.45950955812580774, 0.67787134334312127, 0.70949084487542391, 0.54049044187419226, 0.74002646360041857, 0.25847684349302813)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

$wb.Save()
